# Remove the "Requisitos" section (its Heading2 title paragraph and the
# ListBullet paragraph listing the weak prerequisites that follows it)
# from the end of the document, right after the "Bibliografia" section.

$d = $word.ActiveDocument

# Locate the "Requisitos" heading paragraph by its text + style, so the
# edit does not depend on a brittle, hard-coded paragraph index.
$target = $null
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.Trim()
    if ($txt -eq "Requisitos" -and $p.Style.NameLocal -eq "Heading 2") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $startPos = $target.Range.Start
    $endPos = $d.Content.End

    # Normally the requirements list is the paragraph right after the
    # heading; fall back to the end of the document if, for some reason,
    # the heading turns out to be the very last paragraph.
    $next = $target.Next()
    if ($next -ne $null) {
        $endPos = $next.Range.End
    }

    $r = $d.Range($startPos, $endPos)
    $r.Delete()
}
